# The source diff for this revision touches only word/document.xml and
# word/styles.xml, and every one of its hunks is a pure XML-attribute /
# xmlns-declaration reordering (e.g. `w:tab w:val="left" w:pos="3119"`
# becomes `w:tab w:pos="3119" w:val="left"`, `w:pgSz w:w=".." w:h=".."`
# becomes `w:pgSz w:h=".." w:w=".."`, the root `<w:document>` namespace
# declarations get alphabetized, every `w:lsdException`/`w:style` element
# gets its attributes alphabetized, etc.). Every removed line and its
# paired added line carry the exact same element name and the exact same
# attribute name/value pairs - only the serialized attribute order
# differs. There is no text, formatting, structural, or content change
# anywhere in the diff (matches the commit message: "Fixed POI packaging
# and upgraded to POI 3.15" - upgrading the OOXML library made its XML
# writer alphabetize attributes when it re-serialized this already
# committed test fixture).
#
# The Word object model (what this COM-interop surface exposes) has no
# concept of, or control over, the raw attribute ordering used when the
# underlying XML parts are serialized - that is an internal writer detail
# below the OM. So the only faithful way to reproduce this change through
# Word automation is to leave the document's content/formatting exactly
# as-is: touch nothing, so no semantic diff is introduced.
$d = $word.ActiveDocument
